$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row labels between row 8 and row 9 (shared strings 20/21)
$ws.Range("A8").Value = "econ_program_unitcost_xpert"
$ws.Range("A9").Value = "econ_program_totalcost_xpert"

# Row 8 dummy Xpert cost values
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0
$ws.Range("U8").Value = 0
$ws.Range("V8").Value = 0
$ws.Range("W8").Value = 0
$ws.Range("X8").Value = 0
$ws.Range("Y8").Value = 0
$ws.Range("Z8").Value = 0
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 0
$ws.Range("AC8").Value = 0
$ws.Range("AD8").Value = 0
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 0
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AI8").Value = 0
$ws.Range("AJ8").Value = 0
$ws.Range("AK8").Value = 0
$ws.Range("AL8").Value = 0
$ws.Range("AM8").Value = 0
$ws.Range("AN8").Value = 0
$ws.Range("AO8").Value = 0
$ws.Range("AP8").Value = 0
$ws.Range("AQ8").Value = 0
$ws.Range("AR8").Value = 0
$ws.Range("AS8").Value = 0
$ws.Range("AT8").Value = 0
$ws.Range("AU8").Value = 0
$ws.Range("AV8").Value = 0
$ws.Range("AW8").Value = 0
$ws.Range("AX8").Value = 0
$ws.Range("AY8").Value = 20.77
$ws.Range("BB8").Value = 22
$ws.Range("BE8").Value = 21
$ws.Range("BF8").Value = 21
$ws.Range("BG8").Value = 21
$ws.Range("BH8").Value = 21
$ws.Range("BI8").Value = 21

# Row 9 dummy Xpert cost values
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("Q9").Value = 30000
$ws.Range("V9").Value = 50000
$ws.Range("AF9").Value = 80000
$ws.Range("AK9").Value = 100000
$ws.Range("AP9").Value = 100000
$ws.Range("AZ9").Value = 100000
$ws.Range("BE9").Value = 110000
$ws.Range("BF9").Value = 250000
$ws.Range("BH9").Value = 250000
$ws.Range("BI9").Value = 250000

# Update selection to match saved view
$ws.Range("A8").Select()
